$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "45×88="; New = "56×12=" },
    @{ Old = "62×92="; New = "76×81=" },
    @{ Old = "30×92="; New = "20×82=" },
    @{ Old = "45×21="; New = "53×84=" },
    @{ Old = "67×77="; New = "17×59=" },
    @{ Old = "72×87="; New = "15×68=" },
    @{ Old = "26×36="; New = "53×84=" },
    @{ Old = "13×70="; New = "72×26=" },
    @{ Old = "45×32="; New = "11×90=" },
    @{ Old = "66×28="; New = "40×13=" },
    @{ Old = "63×83="; New = "50×80=" },
    @{ Old = "20×48="; New = "59×58=" },
    @{ Old = "99×62="; New = "20×39=" },
    @{ Old = "84×87="; New = "63×84=" },
    @{ Old = "46×96="; New = "80×74=" },
    @{ Old = "71×83="; New = "33×91=" },
    @{ Old = "24×44="; New = "38×44=" },
    @{ Old = "35×95="; New = "70×55=" },
    @{ Old = "62×91="; New = "92×67=" },
    @{ Old = "16×13="; New = "97×70=" },
    @{ Old = "80×76="; New = "25×99=" },
    @{ Old = "24×34="; New = "22×46=" },
    @{ Old = "26×20="; New = "33×93=" },
    @{ Old = "14×72="; New = "65×54=" },
    @{ Old = "66×74="; New = "92×51=" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.New, 2)
}
